$p = $ppt.ActivePresentation
Write-Output "HasHandoutMaster=$($p.HasHandoutMaster)"
$hm = $p.HandoutMaster
if ($null -eq $hm) {
  Write-Output "HM NULL"
} else {
  Write-Output "HM ok"
  $th = $hm.Theme
  Write-Output "theme: $th"
}
